# Update "想去人数" (want-to-go count) values in column F for rows 2-5
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 41
    $ws.Range("F3").Value = 141
    $ws.Range("F4").Value = 15
    $ws.Range("F5").Value = 35
}
